$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (plain/default, no border/alignment/number-format) taken
# from an existing, unmodified "codice_particella" cell -- used to strip the
# transient Text number-format applied below so purely-numeric-looking codes
# (e.g. "766", "605") are stored as text without leaving a stray cell style.
$plainStyle = $ws.Cells.Item(3, 2).Style

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $plainStyle
}

# A new "particella" (766) was found for comune 140, so it is inserted at row
# 27 and every following record in columns B (codice_particella) and C
# (codice_comune_catastale) shifts down by one row; column A (the row index)
# already lines up and simply gains one more entry (65) at the new last row.
Set-TextCell 27 2 "766"
$ws.Cells.Item(27, 3).Value = 140
Set-TextCell 28 2 ".1512"
$ws.Cells.Item(28, 3).Value = 140
$ws.Cells.Item(29, 2).Value = ".7."
$ws.Cells.Item(29, 3).Value = 187
$ws.Cells.Item(30, 2).Value = "2727/1"
$ws.Cells.Item(30, 3).Value = 189
$ws.Cells.Item(31, 2).Value = "47/3"
$ws.Cells.Item(31, 3).Value = 277
$ws.Cells.Item(32, 2).Value = "302/1"
$ws.Cells.Item(32, 3).Value = 277
$ws.Cells.Item(33, 2).Value = "2103/7"
$ws.Cells.Item(33, 3).Value = 394
$ws.Cells.Item(34, 2).Value = "673/2"
$ws.Cells.Item(34, 3).Value = 384
Set-TextCell 35 2 ".372"
$ws.Cells.Item(35, 3).Value = 384
Set-TextCell 36 2 ".373"
$ws.Cells.Item(36, 3).Value = 384
Set-TextCell 37 2 ".374"
$ws.Cells.Item(37, 3).Value = 384
$ws.Cells.Item(38, 2).Value = "673/2"
$ws.Cells.Item(38, 3).Value = 384
$ws.Cells.Item(39, 2).Value = "406/3"
$ws.Cells.Item(39, 3).Value = 384
Set-TextCell 40 2 "605"
$ws.Cells.Item(40, 3).Value = 384
$ws.Cells.Item(41, 2).Value = "657/1"
$ws.Cells.Item(41, 3).Value = 384
Set-TextCell 42 2 "674"
$ws.Cells.Item(42, 3).Value = 384
$ws.Cells.Item(43, 2).Value = "765/3"
$ws.Cells.Item(43, 3).Value = 384
Set-TextCell 44 2 "938"
$ws.Cells.Item(44, 3).Value = 384
Set-TextCell 45 2 "996"
$ws.Cells.Item(45, 3).Value = 384
Set-TextCell 46 2 "2074"
$ws.Cells.Item(46, 3).Value = 384
Set-TextCell 47 2 "2050"
$ws.Cells.Item(47, 3).Value = 384
Set-TextCell 48 2 "2065"
$ws.Cells.Item(48, 3).Value = 384
Set-TextCell 49 2 "2066"
$ws.Cells.Item(49, 3).Value = 384
Set-TextCell 50 2 "2153"
$ws.Cells.Item(50, 3).Value = 384
Set-TextCell 51 2 "2154"
$ws.Cells.Item(51, 3).Value = 384
Set-TextCell 52 2 "1419"
$ws.Cells.Item(52, 3).Value = 287
Set-TextCell 53 2 "1420"
$ws.Cells.Item(53, 3).Value = 287
Set-TextCell 54 2 "1421"
$ws.Cells.Item(54, 3).Value = 287
Set-TextCell 55 2 "1430"
$ws.Cells.Item(55, 3).Value = 287
$ws.Cells.Item(56, 2).Value = "1431/1"
$ws.Cells.Item(56, 3).Value = 287
$ws.Cells.Item(57, 2).Value = "1431/34"
$ws.Cells.Item(57, 3).Value = 287
$ws.Cells.Item(58, 2).Value = ".950/1"
$ws.Cells.Item(58, 3).Value = 287
$ws.Cells.Item(59, 2).Value = ".950/2"
$ws.Cells.Item(59, 3).Value = 287
Set-TextCell 60 2 "239"
$ws.Cells.Item(60, 3).Value = 287
Set-TextCell 61 2 "241"
$ws.Cells.Item(61, 3).Value = 287
Set-TextCell 62 2 "773"
$ws.Cells.Item(62, 3).Value = 441
Set-TextCell 63 2 ".4046"
$ws.Cells.Item(63, 3).Value = 307
Set-TextCell 64 2 ".4047"
$ws.Cells.Item(64, 3).Value = 307
Set-TextCell 65 2 "1133"
$ws.Cells.Item(65, 3).Value = 248
$ws.Cells.Item(66, 2).Value = "1585/60"
$ws.Cells.Item(66, 3).Value = 248
$ws.Cells.Item(67, 2).Value = "1064/3"
$ws.Cells.Item(67, 3).Value = 248

# Extend the column-A index sequence for the new final row, copying row 66s
# cell formatting (bold, centered, bordered) onto the freshly created A67.
$ws.Cells.Item(67, 1).Value = 65
$ws.Range("A66").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$excel.CutCopyMode = $false
